$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$data = @(
    @(44460, 0, 2, 33.27233405423391),
    @(44461, 0, 2, 33.27233405423391),
    @(44462, 0, 2, 33.27233405423391),
    @(44463, 3, 3, 49.90850108135086),
    @(44464, 0, 3, 49.90850108135086),
    @(44465, 3, 6, 99.81700216270171),
    @(44466, 0, 6, 99.81700216270171),
    @(44467, 0, 6, 99.81700216270171),
    @(44468, 0, 6, 99.81700216270171),
    @(44469, 0, 6, 99.81700216270171),
    @(44470, 0, 3, 49.90850108135086),
    @(44471, 0, 3, 49.90850108135086),
    @(44472, 0, 0, 0),
    @(44473, 0, 0, 0),
    @(44474, 0, 0, 0),
    @(44475, 1, 1, 16.63616702711695),
    @(44476, 1, 2, 33.27233405423391),
    @(44477, 0, 2, 33.27233405423391),
    @(44478, 0, 2, 33.27233405423391),
    @(44479, 0, 2, 33.27233405423391),
    @(44480, 0, 2, 33.27233405423391),
    @(44481, 0, 2, 33.27233405423391),
    @(44482, 0, 1, 16.63616702711695),
    @(44483, 0, 0, 0),
    @(44484, 0, 0, 0),
    @(44485, 0, 0, 0),
    @(44486, 0, 0, 0),
    @(44487, 0, 0, 0),
    @(44488, 0, 0, 0),
    @(44489, 0, 0, 0),
    @(44490, 0, 0, 0),
    @(44491, 0, 0, 0),
    @(44492, 0, 0, 0),
    @(44493, 0, 0, 0),
    @(44494, 0, 0, 0),
    @(44495, 0, 0, 0),
    @(44496, 0, 0, 0),
    @(44497, 0, 0, 0),
    @(44498, 1, 1, 16.63616702711695),
    @(44499, 0, 1, 16.63616702711695),
    @(44500, 2, 3, 49.90850108135086),
    @(44501, 8, 11, 182.9978372982865),
    @(44502, 3, 14, 232.9063383796373),
    @(44503, 0, 14, 232.9063383796373),
    @(44504, 0, 14, 232.9063383796373),
    @(44505, 9, 22, 365.9956745965729),
    @(44506, 2, 24, 399.2680086508068),
    @(44507, 4, 26, 432.5403427050408),
    @(44508, 3, 21, 349.359507569456),
    @(44509, 3, 21, 349.359507569456),
    @(44510, 0, 21, 349.359507569456),
    @(44511, 1, 22, 365.9956745965729),
    @(44512, 6, 19, 316.0871735152221),
    @(44513, 3, 20, 332.723340542339),
    @(44514, 0, 16, 266.1786724338713),
    @(44515, 0, 13, 216.2701713525204),
    @(44516, 2, 12, 199.6340043254034),
    @(44517, 0, 12, 199.6340043254034),
    @(44518, 0, 11, 182.9978372982865),
    @(44519, 4, 9, 149.7255032440526),
    @(44520, 0, 6, 99.81700216270171),
    @(44521, 1, 7, 116.4531691898187),
    @(44522, 2, 9, 149.7255032440526),
    @(44523, 1, 8, 133.0893362169356),
    @(44524, 2, 10, 166.3616702711695),
    @(44525, 2, 12, 199.6340043254034),
    @(44526, 2, 10, 166.3616702711695),
    @(44527, 0, 10, 166.3616702711695),
    @(44528, 1, 10, 166.3616702711695),
    @(44529, 2, 10, 166.3616702711695),
    @(44530, 6, 15, 249.5425054067543),
    @(44531, 0, 13, 216.2701713525204),
    @(44532, 0, 11, 182.9978372982865),
    @(44533, 3, 12, 199.6340043254034),
    @(44534, 0, 12, 199.6340043254034),
    @(44535, 4, 15, 249.5425054067543),
    @(44536, 3, 16, 266.1786724338713),
    @(44537, 3, 13, 216.2701713525204),
    @(44538, 0, 13, 216.2701713525204),
)

$startRow = 386
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

# Copy the date-column style (xf index 2: centered, bordered, bold, yyyy-mm-dd format)
# from the last pre-existing data row down onto the newly added date cells.
$endRow = $startRow + $data.Count - 1
$ws.Range("A385").Copy()
$ws.Range("A" + $startRow + ":A" + $endRow).PasteSpecial(-4122)
$excel.CutCopyMode = $false
